$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("input")
$wsExplanation = $wb.Worksheets.Item("explanation")

# --- Update the "input" sheet data row (row 2) ---
# Analysis name changed from test1 -> Nepal
$wsInput.Range("A2").Value = "Nepal"
# Network source changed from "Network based on shapefile" -> "Network based on OSM online"
$wsInput.Range("D2").Value = "Network based on OSM online"
# New OSM area of interest shapefile name (previously blank)
$wsInput.Range("E2").Value = "npl_admbnda_adm0_nd_20201117.shp"
# Old shapefile-specific columns (part_of_DR_roads / fid) no longer apply -> clear them
$wsInput.Range("G2").ClearContents() | Out-Null
$wsInput.Range("H2").ClearContents() | Out-Null
# Fill in network_type / road_types that used to be blank
$wsInput.Range("L2").Value = "drive"
$wsInput.Range("M2").Value = "motorway, trunk, primary, secondary"

# --- Update sheet selections / active tab ---
# Previously "explanation" was the active/selected tab with selection E2;
# make "input" the active tab with the whole second row selected, and
# leave "explanation" selected at A2 (no longer the active tab).
$wsExplanation.Range("A2").Select() | Out-Null
$wsInput.Activate()
$wsInput.Range("A2:XFD2").Select() | Out-Null
